$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.786.49"
$ws.Range("E2").Value = "  +5.79%  "
$ws.Range("D3").Value = "3.111.34"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.29"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.29"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.101.14"
$ws.Range("E8").Value = "  +3.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("E10").Value = "  +8.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("E11").Value = "  +10.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").Value = "  +4.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.56"
$ws.Range("E14").Value = "  +5.00%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.625.05"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.28"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.109.78"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "62.741.29"
$ws.Range("E19").Value = "  +5.78%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "452.94"
$ws.Range("E20").Value = "  +4.60%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.09"
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.735"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  +5.26%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.66"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.07"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.25"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("E28").Value = "  +5.44%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.26"
$ws.Range("E30").Value = "  +4.59%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.84"
$ws.Range("E31").Value = "  +11.45%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.112"
$ws.Range("E32").Value = "  +12.14%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.14"
$ws.Range("E33").Value = "  +4.82%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.05"
$ws.Range("E34").Value = "  +4.64%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0803"
$ws.Range("E35").Value = "  +4.74%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.08"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.22"
$ws.Range("E37").Value = "  +4.72%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.58"
$ws.Range("E38").Value = "  +3.40%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.01"
$ws.Range("E39").Value = "  +9.12%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.80"
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "425.17"
$ws.Range("E41").Value = "  +4.76%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.945.36"
$ws.Range("E42").Value = "  +6.19%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0374"
$ws.Range("E43").Value = "  +5.24%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.279"
$ws.Range("E44").Value = "  +9.72%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.111"
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  +6.43%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.69"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.10"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.112"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.86"
$ws.Range("E51").Value = "  +5.27%  "
